# Generate Report for Handoff
# - Mark the "Ready for handoff" rows (8, 10-14) with Priority "ht" on the
#   per-language handoff sheets.
# - Refresh the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#   timestamps for the same rows on every sheet that surfaces them
#   (zh-cn/de-de "Latest Handoff Datetime" plus the Overview roll-up
#   column, which tracks the de-de timestamp).

$wb = $excel.ActiveWorkbook

$rows = @(8, 10, 11, 12, 13, 14)

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $ws_zhcn.Range("E$r").Value = "ht"
    $ws_zhcn.Range("H$r").Value = "2016-09-03 10:23:46"
}

$ws_dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $ws_dede.Range("E$r").Value = "ht"
    $ws_dede.Range("H$r").Value = "2016-09-03 10:23:52"
}

# Overview sheet mirrors the de-de "Latest Handoff Datetime" column (shared
# text) under "Latest HO Xliff Generate Date".
$ws_overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $ws_overview.Range("G$r").Value = "2016-09-03 10:23:52"
}
